$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.78
$ws.Range("J2").Value = 4.7
$ws.Range("L2").Value = 1.22
$ws.Range("N2").Value = 7.8
$ws.Range("O2").Value = 1.13
$ws.Range("P2").Value = 3.3
$ws.Range("Q2").Value = 1.4
$ws.Range("R2").Value = 1.96
$ws.Range("S2").Value = 1.96
$ws.Range("T2").Value = 1.47
$ws.Range("U2").Value = 2.96
$ws.Range("W2").Value = 2.28
$ws.Range("X2").Value = 34
$ws.Range("Y2").Value = 980
$ws.Range("Z2").Value = 980
$ws.Range("AA2").Value = 100
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 15
$ws.Range("AE2").Value = 980
$ws.Range("AF2").Value = 19.5
$ws.Range("AG2").Value = 13
$ws.Range("AH2").Value = 18
$ws.Range("AI2").Value = 980
$ws.Range("AJ2").Value = 980
$ws.Range("AK2").Value = 18.5
$ws.Range("AL2").Value = 25
$ws.Range("AN2").Value = 5.6
$ws.Range("AO2").Value = 980
$ws.Range("F3").Value = 2.78
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 2.72
$ws.Range("I3").Value = 3.5
$ws.Range("J3").Value = 2.42
$ws.Range("K3").Value = 3.05
$ws.Range("L3").Value = 1.58
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 2.16
$ws.Range("O3").Value = 1.7
$ws.Range("P3").Value = 1.38
$ws.Range("Q3").Value = 2.8
$ws.Range("R3").Value = 1.13
$ws.Range("S3").Value = 1.02
$ws.Range("T3").Value = 2.32
$ws.Range("U3").Value = 1.61
$ws.Range("V3").Value = 1.4
$ws.Range("W3").Value = 1.33
$ws.Range("X3").Value = 8
$ws.Range("Y3").Value = 8.800000000000001
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 70
$ws.Range("AB3").Value = 9.6
$ws.Range("AC3").Value = 8.4
$ws.Range("AD3").Value = 18
$ws.Range("AE3").Value = 65
$ws.Range("AF3").Value = 25
$ws.Range("AG3").Value = 20
$ws.Range("AH3").Value = 36
$ws.Range("AI3").Value = 120
$ws.Range("AJ3").Value = 90
$ws.Range("AK3").Value = 80
$ws.Range("AL3").Value = 130
$ws.Range("AM3").Value = 350
$ws.Range("AN3").Value = 120
$ws.Range("AO3").Value = 90
$ws.Range("F4").Value = 4.5
$ws.Range("G4").Value = 6.6
$ws.Range("H4").Value = 1.78
$ws.Range("I4").Value = 1.94
$ws.Range("J4").Value = 3.35
$ws.Range("L4").Value = 1.41
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 3.05
$ws.Range("P4").Value = 1.71
$ws.Range("Q4").Value = 2.14
$ws.Range("R4").Value = 1.26
$ws.Range("S4").Value = 3.75
$ws.Range("T4").Value = 1.98
$ws.Range("V4").Value = 2.06
$ws.Range("W4").Value = 1.18
$ws.Range("Y4").Value = 8.4
$ws.Range("Z4").Value = 12
$ws.Range("AA4").Value = 980
$ws.Range("AC4").Value = 9.4
$ws.Range("AE4").Value = 980
$ws.Range("AF4").Value = 980
$ws.Range("AG4").Value = 980
$ws.Range("AH4").Value = 980
$ws.Range("AI4").Value = 980
$ws.Range("AK4").Value = 100
$ws.Range("AL4").Value = 110
$ws.Range("AN4").Value = 140
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 1.34
$ws.Range("N5").Value = 5.7
$ws.Range("O5").Value = 1.19
$ws.Range("P5").Value = 2.62
$ws.Range("Q5").Value = 1.56
$ws.Range("R5").Value = 1.64
$ws.Range("S5").Value = 2.42
$ws.Range("T5").Value = 1.94
$ws.Range("X5").Value = 28
$ws.Range("Y5").Value = 11.5
$ws.Range("Z5").Value = 9.199999999999999
$ws.Range("AA5").Value = 11.5
$ws.Range("AB5").Value = 40
$ws.Range("AC5").Value = 14
$ws.Range("AD5").Value = 11.5
$ws.Range("AE5").Value = 15
$ws.Range("AF5").Value = 110
$ws.Range("AG5").Value = 40
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 36
$ws.Range("AJ5").Value = 400
$ws.Range("AK5").Value = 170
$ws.Range("AL5").Value = 140
$ws.Range("AM5").Value = 150
$ws.Range("AN5").Value = 200
$ws.Range("AO5").Value = 4.7
$ws.Range("F6").Value = 2.7
$ws.Range("G6").Value = 2.76
$ws.Range("O6").Value = 1.45
$ws.Range("P6").Value = 1.7
$ws.Range("Q6").Value = 2.34
$ws.Range("T6").Value = 1.97
$ws.Range("U6").Value = 1.97
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 9.800000000000001
$ws.Range("Z6").Value = 19.5
$ws.Range("AA6").Value = 55
$ws.Range("AB6").Value = 9.199999999999999
$ws.Range("AC6").Value = 7.2
$ws.Range("AD6").Value = 13.5
$ws.Range("AE6").Value = 40
$ws.Range("AF6").Value = 16.5
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 20
$ws.Range("AI6").Value = 60
$ws.Range("AJ6").Value = 44
$ws.Range("AK6").Value = 36
$ws.Range("AL6").Value = 55
$ws.Range("AM6").Value = 160
$ws.Range("AN6").Value = 36
$ws.Range("AO6").Value = 44
$ws.Range("H7").Value = 21
$ws.Range("P7").Value = 2.6
$ws.Range("S7").Value = 2.48
$ws.Range("T7").Value = 2.62
$ws.Range("X7").Value = 29
$ws.Range("Y7").Value = 60
$ws.Range("AB7").Value = 9.199999999999999
$ws.Range("AC7").Value = 20
$ws.Range("AD7").Value = 85
$ws.Range("AF7").Value = 7
$ws.Range("AG7").Value = 13
$ws.Range("AH7").Value = 55
$ws.Range("AI7").Value = 410
$ws.Range("AJ7").Value = 7.8
$ws.Range("AK7").Value = 15.5
$ws.Range("AL7").Value = 60
$ws.Range("AM7").Value = 440
$ws.Range("AN7").Value = 3.7
$ws.Range("H8").Value = 2.96
$ws.Range("J8").Value = 3.4
$ws.Range("M8").Value = 1.08
$ws.Range("X8").Value = 13.5
$ws.Range("Y8").Value = 11.5
$ws.Range("Z8").Value = 20
$ws.Range("AA8").Value = 50
$ws.Range("AB8").Value = 11
$ws.Range("AD8").Value = 13
$ws.Range("AE8").Value = 36
$ws.Range("AF8").Value = 17.5
$ws.Range("AG8").Value = 12.5
$ws.Range("AH8").Value = 17.5
$ws.Range("AI8").Value = 46
$ws.Range("AJ8").Value = 40
$ws.Range("AK8").Value = 30
$ws.Range("AL8").Value = 44
$ws.Range("AM8").Value = 120
$ws.Range("AN8").Value = 26
$ws.Range("AO8").Value = 32
$ws.Range("G9").Value = 2.66
$ws.Range("H9").Value = 3.65
$ws.Range("K9").Value = 3
$ws.Range("F10").Value = 2.28
$ws.Range("H10").Value = 3.95
$ws.Range("I10").Value = 4.4
$ws.Range("K10").Value = 3.15
